# Updated cryptos list row data pulled from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D="68.574.17"; E="  +0.03%  " },
    @{ Row=3; D="3.906.01"; E="  +0.08%  " },
    @{ Row=4; E="  -0.08%  " },
    @{ Row=5; D="603.18"; E="  +0.09%  " },
    @{ Row=6; D="169.66"; E="  +1.62%  " },
    @{ Row=7; D="3.904.34"; E="  +0.06%  " },
    @{ Row=8; E="  +0.05%  " },
    @{ Row=9; E="  +0.67%  " },
    @{ Row=10; E="  -0.52%  " },
    @{ Row=11; D="6.45"; E="  +0.05%  " },
    @{ Row=12; E="  +0.18%  " },
    @{ Row=13; D="0.0000255"; E="  -0.16%  " },
    @{ Row=14; D="37.31"; E="  -0.34%  " },
    @{ Row=15; D="4.564.13"; E="  +0.16%  " },
    @{ Row=16; D="3.911.35"; E="  +0.30%  " },
    @{ Row=17; D="68.540.54"; E="  -0.22%  " },
    @{ Row=18; D="18.18"; E="  +4.98%  " },
    @{ Row=19; D="7.45"; E="  -0.22%  " },
    @{ Row=20; E="  +0.37%  " },
    @{ Row=21; D="10.86"; E="  -1.74%  " },
    @{ Row=22; D="473.40"; E="  -3.09%  " },
    @{ Row=23; D="0.744"; E="  +2.31%  " },
    @{ Row=24; E="  +0.91%  " },
    @{ Row=25; D="83.92"; E="  -0.89%  " },
    @{ Row=26; E="  +1.17%  " },
    @{ Row=27; D="12.26"; E="  +1.77%  " },
    @{ Row=28; D="10.05"; E="  -1.07%  " },
    @{ Row=29; E="  +0.16%  " },
    @{ Row=31; D="4.056.12"; E="  +0.02%  " },
    @{ Row=32; D="7.91"; E="  +2.32%  " },
    @{ Row=33; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="2.33"; E="  -2.13%  " },
    @{ Row=34; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="31.56"; E="  -0.90%  " },
    @{ Row=35; D="9.50"; E="  +1.88%  " },
    @{ Row=36; D="3.877.89"; E="  +0.52%  " },
    @{ Row=37; E="  -1.11%  " },
    @{ Row=38; D="3.72"; E="  +16.01%  " },
    @{ Row=39; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.142"; E="  +2.02%  " },
    @{ Row=40; B="Mantle"; C="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D="1.03"; E="  -0.48%  " },
    @{ Row=41; E="  +0.22%  " },
    @{ Row=42; D="1.00"; E="  +0.00%  " },
    @{ Row=43; E="  -0.62%  " },
    @{ Row=44; B="FLOKI"; C="https://coinranking.com/coin/fmHk13Rqw+floki-floki"; D="0.000303"; E="  +13.78%  " },
    @{ Row=45; B="Bittensor"; C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D="429.91"; E="  +0.08%  " },
    @{ Row=46; B="Stacks"; C="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D="2.00"; E="  +0.56%  " },
    @{ Row=47; D="8.68"; E="  +1.76%  " },
    @{ Row=48; E="  +0.03%  " },
    @{ Row=49; D="47.16"; E="  -1.98%  " },
    @{ Row=50; D="27.11"; E="  +5.96%  " },
    @{ Row=51; D="144.04"; E="  +0.75%  " }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($r.ContainsKey("B")) { $ws.Cells.Item($rowNum, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($rowNum, 3).Value = $r.C }
    if ($r.ContainsKey("D")) {
        $cell = $ws.Cells.Item($rowNum, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $r.D
    }
    if ($r.ContainsKey("E")) { $ws.Cells.Item($rowNum, 5).Value = $r.E }
}
